# Reorders the "Requisitos" bullet list in LOB1237.docx so the individual
# requirement lines appear in the new order recorded by the site build
# (same 27 requirement lines -- pure reshuffle, no text/formatting changes).

$d = $word.ActiveDocument

# Locate the paragraph that holds the requirement list. It is the
# ListBullet paragraph whose text begins with the first requirement line
# ("LOB1268 - ..."). Using -like with a trailing "*" because Range.Text
# includes the trailing paragraph-mark character.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "LOB1268*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Host "ERROR: could not locate the Requisitos list paragraph"
} else {
    $startPos = $target.Range.Start
    $endPos = $target.Range.End

    # Build a fresh Range over the whole paragraph content (excluding the
    # paragraph mark) and overwrite it via InsertXML so every requirement
    # keeps living in its own <w:r>...<w:br/></w:r> run, matching the
    # original document's structure -- only the order of the runs changes.
    $rng = $d.Range($startPos, $endPos)

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>LOB1268 -  Leitura, Escrita e Comunicação Científica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1270 -  Química Experimental Aplicada  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4233 -  Gestão de Negócios  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1011 -  Eletricidade Aplicada  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1052 -  Cálculo III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1019 -  Física II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1038 -  Física Experimental I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1053 -  Física III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1037 -  Álgebra Linear  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1042 -  Física Experimental IV  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1003 -  Cálculo I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1041 -  Física Experimental II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1018 -  Física I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1006 -  Cálculo IV  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1039 -  Física Experimental III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1021 -  Física IV  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1036 -  Geometria Analítica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1024 -  Mecânica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1040 -  Laboratório de Eletricidade  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1004 -  Cálculo II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1012 -  Estatística  (Requisito)</w:t><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xml)

    Write-Host "Requisitos list reordered successfully."
}
